$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.11"
$meta.Range("B8").Value = "2024-06-13T17:23:26-04:00"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Row 6 = Extension.value[x] (base, unsliced) -- slicing becomes "open" and the
# Type(s) column now lists every permitted data type instead of just "Coding".
$typesList = "base64Binary`nbooleancanonicalcodedatedateTimedecimalidinstantintegermarkdownoidpositiveIntstringtimeunsignedInturiurluuidAddressAgeAnnotationAttachmentCodeableConceptCodingContactPointCountDistanceDurationHumanNameIdentifierMoneyPeriodQuantityRangeRatioReferenceSampledDataSignatureTimingContactDetailContributorDataRequirementExpressionParameterDefinitionRelatedArtifactTriggerDefinitionUsageContextDosageMeta"
$elements.Range("K6").Value = $typesList
$elements.Range("AE6").Value = "open"

# Widen column K (Type(s)) to fit the newly-expanded list of types.
$elements.Columns("K").ColumnWidth = 254.16666666666666
